$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "x"
$ws.Range("A3").Value = "X"
$ws.Range("A4").Value = "X"
$ws.Range("A5").Value = "X"

$ws.Range("B2").Value = "o"
$ws.Range("B3").Value = "o"
$ws.Range("B4").Value = "o"
$ws.Range("B5").Value = "o"

$ws.Range("E6").Select()
